# Crear/Actualizar Excel para pedido 6903661db7cb420aeabeeab5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: D12 becomes empty (remove its content entirely)
$ws.Range("D12").ClearContents()

# Row 13: new record
$ws.Cells.Item(13, 1).Value = 2187
$ws.Cells.Item(13, 2).Value = "Samuel Rupérez Macarro"
$ws.Cells.Item(13, 3).Value = "Estructura coplanar NOVOTEGRA"
$ws.Cells.Item(13, 4).Value = "MODULO FV JA SOLAR 535WP BLACK FRAME BIFACIAL 120 CELDAS"

$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "16"
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(13, 8).Value = "GOODWE ES UNIQ - GW8000-ES-C10 híbrido monofásico"

$ws.Cells.Item(13, 9).NumberFormat = "@"
$ws.Cells.Item(13, 9).Value = "1"
$ws.Cells.Item(13, 9).Style = "Normal"

$ws.Cells.Item(13, 10).Value = "GOODWE Batería Lynx Home U G3 5,12 kWh"

$ws.Cells.Item(13, 11).NumberFormat = "@"
$ws.Cells.Item(13, 11).Value = "1"
$ws.Cells.Item(13, 11).Style = "Normal"

$ws.Cells.Item(13, 13).Value = "Sí"
$ws.Cells.Item(13, 14).Value = "2025-09-25T07:50:43.054Z"
